$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.878.19'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.843.11'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'309.27"
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = '  +3.44%  '
$ws.Range("D8").Value = "'0.3659"
$ws.Range("E8").Value = '  +1.94%  '
$ws.Range("D9").Value = "'0.07150"
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = "'0.9260"
$ws.Range("E10").Value = '  +3.79%  '
$ws.Range("D11").Value = "'19.56"
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = "'0.07685"
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '1.869.57'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").Value = "'6.393"
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = "'88.28"
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = "'0.000008642"
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '26.925.49'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").Value = "'5.014"
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("D24").Value = "'1.931"
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").Value = "'151.78"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'18.23"
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").Value = "'114.23"
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = "'4.876"
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("D30").Value = "'0.08833"
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").Value = "'1.177"
$ws.Range("E32").Value = '  +6.07%  '
$ws.Range("D33").Value = "'0.7483"
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").Value = "'2.775"
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("D35").Value = "'4.469"
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").Value = "'1.082"
$ws.Range("D37").Value = "'0.01940"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = "'0.05206"
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("D39").Value = "'2.952"
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("D40").Value = "'0.5217"
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").Value = "'6.982"
$ws.Range("E41").Value = '  +3.12%  '
$ws.Range("D42").Value = "'0.1512"
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").Value = "'8.149"
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("E44").Value = '  +5.81%  '
$ws.Range("D45").Value = "'0.4700"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = "'1.005"
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").Value = "'100.91"
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").Value = "'1.596"
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").Value = "'65.46"
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("D51").Value = "'0.8918"
$ws.Range("E51").Value = '  +5.28%  '
